$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update aggregate counts that shifted because the "Goa, Índia" (row 48)
# entry was reclassified/removed from the breakdown.
$ws.Range("B2").Value = 506
$ws.Range("B3").Value = 47
$ws.Range("B7").Value = 31

# Remove the "Goa, Índia" row entirely; everything below shifts up one row
# (Excel's row delete also keeps the used-range dimension in sync).
$ws.Rows.Item(48).Delete()
